$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$t = $s.Shapes.Item(2).Table
$t.ApplyStyle("{9A91CC96-4B0C-478E-A8C5-8F830C9C961A}")
